# ===========================================================================
# Add 2022-Q1 fund-holdings data:
#   1. Insert a new "2022-Q1" worksheet (positioned right before "总计"),
#      populated with the quarter's fund holdings table.
#   2. Prepend a "2022-Q1" summary row to the "总计" worksheet, shifting the
#      existing rows down and renumbering the index column.
# ===========================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "2022-Q1" sheet
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($total)
$ws.Name = "2022-Q1"

# Header row (B1:H1)
$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"

# Data rows (A2:H38) - fund code/size/position columns are stored as text
# (leading "'" forces text so numeric-looking strings like "30.22" keep
# their original text representation instead of becoming Number cells).
    # row 2
    $ws.Cells.Item(2, 1).Value = 0
    $ws.Cells.Item(2, 2).Value = "'340008"
    $ws.Cells.Item(2, 3).Value = "兴全有机增长混合"
    $ws.Cells.Item(2, 4).Value = "'30.22"
    $ws.Cells.Item(2, 5).Value = "'75.70"
    $ws.Cells.Item(2, 6).Value = "'9.43"
    $ws.Cells.Item(2, 7).Value = "'2.8497"
    $ws.Cells.Item(2, 8).Value = 1
    # row 3
    $ws.Cells.Item(3, 1).Value = 1
    $ws.Cells.Item(3, 2).Value = "'001576"
    $ws.Cells.Item(3, 3).Value = "国泰智能装备股票A"
    $ws.Cells.Item(3, 4).Value = "'44.49"
    $ws.Cells.Item(3, 5).Value = "'91.61"
    $ws.Cells.Item(3, 6).Value = "'4.45"
    $ws.Cells.Item(3, 7).Value = "'1.9798"
    $ws.Cells.Item(3, 8).Value = 9
    # row 4
    $ws.Cells.Item(4, 1).Value = 2
    $ws.Cells.Item(4, 2).Value = "'012748"
    $ws.Cells.Item(4, 3).Value = "华泰柏瑞远见智选混合型证券投资基金A"
    $ws.Cells.Item(4, 4).Value = "'26.44"
    $ws.Cells.Item(4, 5).Value = "'94.15"
    $ws.Cells.Item(4, 6).Value = "'7.13"
    $ws.Cells.Item(4, 7).Value = "'1.8852"
    $ws.Cells.Item(4, 8).Value = 5
    # row 5
    $ws.Cells.Item(5, 1).Value = 3
    $ws.Cells.Item(5, 2).Value = "'011037"
    $ws.Cells.Item(5, 3).Value = "富国长期成长混合型证券投资基金"
    $ws.Cells.Item(5, 4).Value = "'57.98"
    $ws.Cells.Item(5, 5).Value = "'75.37"
    $ws.Cells.Item(5, 6).Value = "'2.14"
    $ws.Cells.Item(5, 7).Value = "'1.2408"
    $ws.Cells.Item(5, 8).Value = 8
    # row 6
    $ws.Cells.Item(6, 1).Value = 4
    $ws.Cells.Item(6, 2).Value = "'100022"
    $ws.Cells.Item(6, 3).Value = "富国天瑞强势地区精选混合"
    $ws.Cells.Item(6, 4).Value = "'53.79"
    $ws.Cells.Item(6, 5).Value = "'75.81"
    $ws.Cells.Item(6, 6).Value = "'2.07"
    $ws.Cells.Item(6, 7).Value = "'1.1135"
    $ws.Cells.Item(6, 8).Value = 9
    # row 7
    $ws.Cells.Item(7, 1).Value = 5
    $ws.Cells.Item(7, 2).Value = "'460001"
    $ws.Cells.Item(7, 3).Value = "华泰柏瑞盛世中国混合"
    $ws.Cells.Item(7, 4).Value = "'18.38"
    $ws.Cells.Item(7, 5).Value = "'85.97"
    $ws.Cells.Item(7, 6).Value = "'5.78"
    $ws.Cells.Item(7, 7).Value = "'1.0624"
    $ws.Cells.Item(7, 8).Value = 6
    # row 8
    $ws.Cells.Item(8, 1).Value = 6
    $ws.Cells.Item(8, 2).Value = "'011322"
    $ws.Cells.Item(8, 3).Value = "国泰智能装备股票C"
    $ws.Cells.Item(8, 4).Value = "'21.46"
    $ws.Cells.Item(8, 5).Value = "'91.61"
    $ws.Cells.Item(8, 6).Value = "'4.45"
    $ws.Cells.Item(8, 7).Value = "'0.9550"
    $ws.Cells.Item(8, 8).Value = 9
    # row 9
    $ws.Cells.Item(9, 1).Value = 7
    $ws.Cells.Item(9, 2).Value = "'009092"
    $ws.Cells.Item(9, 3).Value = "富国新材料新能源混合"
    $ws.Cells.Item(9, 4).Value = "'14.70"
    $ws.Cells.Item(9, 5).Value = "'89.49"
    $ws.Cells.Item(9, 6).Value = "'4.78"
    $ws.Cells.Item(9, 7).Value = "'0.7027"
    $ws.Cells.Item(9, 8).Value = 6
    # row 10
    $ws.Cells.Item(10, 1).Value = 8
    $ws.Cells.Item(10, 2).Value = "'010330"
    $ws.Cells.Item(10, 3).Value = "东吴兴享成长混合A"
    $ws.Cells.Item(10, 4).Value = "'11.63"
    $ws.Cells.Item(10, 5).Value = "'80.15"
    $ws.Cells.Item(10, 6).Value = "'5.73"
    $ws.Cells.Item(10, 7).Value = "'0.6664"
    $ws.Cells.Item(10, 8).Value = 2
    # row 11
    $ws.Cells.Item(11, 1).Value = 9
    $ws.Cells.Item(11, 2).Value = "'180003"
    $ws.Cells.Item(11, 3).Value = "银华-道琼斯88指数 A"
    $ws.Cells.Item(11, 4).Value = "'14.84"
    $ws.Cells.Item(11, 5).Value = "'83.54"
    $ws.Cells.Item(11, 6).Value = "'2.78"
    $ws.Cells.Item(11, 7).Value = "'0.4126"
    $ws.Cells.Item(11, 8).Value = 10
    # row 12
    $ws.Cells.Item(12, 1).Value = 10
    $ws.Cells.Item(12, 2).Value = "'519110"
    $ws.Cells.Item(12, 3).Value = "浦银安盛价值成长混合A"
    $ws.Cells.Item(12, 4).Value = "'8.04"
    $ws.Cells.Item(12, 5).Value = "'88.87"
    $ws.Cells.Item(12, 6).Value = "'3.01"
    $ws.Cells.Item(12, 7).Value = "'0.2420"
    $ws.Cells.Item(12, 8).Value = 3
    # row 13
    $ws.Cells.Item(13, 1).Value = 11
    $ws.Cells.Item(13, 2).Value = "'007163"
    $ws.Cells.Item(13, 3).Value = "浦银安盛环保新能源混合A"
    $ws.Cells.Item(13, 4).Value = "'4.70"
    $ws.Cells.Item(13, 5).Value = "'79.84"
    $ws.Cells.Item(13, 6).Value = "'5.12"
    $ws.Cells.Item(13, 7).Value = "'0.2406"
    $ws.Cells.Item(13, 8).Value = 4
    # row 14
    $ws.Cells.Item(14, 1).Value = 12
    $ws.Cells.Item(14, 2).Value = "'519170"
    $ws.Cells.Item(14, 3).Value = "浦银安盛增长动力灵活配置混合"
    $ws.Cells.Item(14, 4).Value = "'8.12"
    $ws.Cells.Item(14, 5).Value = "'85.61"
    $ws.Cells.Item(14, 6).Value = "'2.92"
    $ws.Cells.Item(14, 7).Value = "'0.2371"
    $ws.Cells.Item(14, 8).Value = 6
    # row 15
    $ws.Cells.Item(15, 1).Value = 13
    $ws.Cells.Item(15, 2).Value = "'010345"
    $ws.Cells.Item(15, 3).Value = "华泰柏瑞成长智选混合A"
    $ws.Cells.Item(15, 4).Value = "'4.70"
    $ws.Cells.Item(15, 5).Value = "'93.47"
    $ws.Cells.Item(15, 6).Value = "'4.71"
    $ws.Cells.Item(15, 7).Value = "'0.2214"
    $ws.Cells.Item(15, 8).Value = 8
    # row 16
    $ws.Cells.Item(16, 1).Value = 14
    $ws.Cells.Item(16, 2).Value = "'012749"
    $ws.Cells.Item(16, 3).Value = "华泰柏瑞远见智选混合型证券投资基金C"
    $ws.Cells.Item(16, 4).Value = "'3.08"
    $ws.Cells.Item(16, 5).Value = "'94.15"
    $ws.Cells.Item(16, 6).Value = "'7.13"
    $ws.Cells.Item(16, 7).Value = "'0.2196"
    $ws.Cells.Item(16, 8).Value = 5
    # row 17
    $ws.Cells.Item(17, 1).Value = 15
    $ws.Cells.Item(17, 2).Value = "'000264"
    $ws.Cells.Item(17, 3).Value = "博时内需增长混合"
    $ws.Cells.Item(17, 4).Value = "'3.87"
    $ws.Cells.Item(17, 5).Value = "'75.26"
    $ws.Cells.Item(17, 6).Value = "'5.29"
    $ws.Cells.Item(17, 7).Value = "'0.2047"
    $ws.Cells.Item(17, 8).Value = 4
    # row 18
    $ws.Cells.Item(18, 1).Value = 16
    $ws.Cells.Item(18, 2).Value = "'166011"
    $ws.Cells.Item(18, 3).Value = "中欧盛世成长混合 (LOF) -A"
    $ws.Cells.Item(18, 4).Value = "'5.80"
    $ws.Cells.Item(18, 5).Value = "'85.98"
    $ws.Cells.Item(18, 6).Value = "'3.50"
    $ws.Cells.Item(18, 7).Value = "'0.2030"
    $ws.Cells.Item(18, 8).Value = 4
    # row 19
    $ws.Cells.Item(19, 1).Value = 17
    $ws.Cells.Item(19, 2).Value = "'001888"
    $ws.Cells.Item(19, 3).Value = "中欧盛世成长混合 (LOF) -E"
    $ws.Cells.Item(19, 4).Value = "'5.80"
    $ws.Cells.Item(19, 5).Value = "'85.98"
    $ws.Cells.Item(19, 6).Value = "'3.50"
    $ws.Cells.Item(19, 7).Value = "'0.2030"
    $ws.Cells.Item(19, 8).Value = 4
    # row 20
    $ws.Cells.Item(20, 1).Value = 18
    $ws.Cells.Item(20, 2).Value = "'007164"
    $ws.Cells.Item(20, 3).Value = "浦银安盛环保新能源混合C"
    $ws.Cells.Item(20, 4).Value = "'3.32"
    $ws.Cells.Item(20, 5).Value = "'79.84"
    $ws.Cells.Item(20, 6).Value = "'5.12"
    $ws.Cells.Item(20, 7).Value = "'0.1700"
    $ws.Cells.Item(20, 8).Value = 4
    # row 21
    $ws.Cells.Item(21, 1).Value = 19
    $ws.Cells.Item(21, 2).Value = "'001306"
    $ws.Cells.Item(21, 3).Value = "中欧永裕混合A"
    $ws.Cells.Item(21, 4).Value = "'4.48"
    $ws.Cells.Item(21, 5).Value = "'86.33"
    $ws.Cells.Item(21, 6).Value = "'3.52"
    $ws.Cells.Item(21, 7).Value = "'0.1577"
    $ws.Cells.Item(21, 8).Value = 4
    # row 22
    $ws.Cells.Item(22, 1).Value = 20
    $ws.Cells.Item(22, 2).Value = "'007306"
    $ws.Cells.Item(22, 3).Value = "华泰柏瑞基本面智选混合A"
    $ws.Cells.Item(22, 4).Value = "'3.54"
    $ws.Cells.Item(22, 5).Value = "'93.61"
    $ws.Cells.Item(22, 6).Value = "'4.16"
    $ws.Cells.Item(22, 7).Value = "'0.1473"
    $ws.Cells.Item(22, 8).Value = 10
    # row 23
    $ws.Cells.Item(23, 1).Value = 21
    $ws.Cells.Item(23, 2).Value = "'050012"
    $ws.Cells.Item(23, 3).Value = "博时策略混合"
    $ws.Cells.Item(23, 4).Value = "'3.00"
    $ws.Cells.Item(23, 5).Value = "'73.86"
    $ws.Cells.Item(23, 6).Value = "'4.49"
    $ws.Cells.Item(23, 7).Value = "'0.1347"
    $ws.Cells.Item(23, 8).Value = 7
    # row 24
    $ws.Cells.Item(24, 1).Value = 22
    $ws.Cells.Item(24, 2).Value = "'519120"
    $ws.Cells.Item(24, 3).Value = "浦银安盛新兴产业混合"
    $ws.Cells.Item(24, 4).Value = "'2.21"
    $ws.Cells.Item(24, 5).Value = "'90.11"
    $ws.Cells.Item(24, 6).Value = "'2.84"
    $ws.Cells.Item(24, 7).Value = "'0.0628"
    $ws.Cells.Item(24, 8).Value = 5
    # row 25
    $ws.Cells.Item(25, 1).Value = 23
    $ws.Cells.Item(25, 2).Value = "'519113"
    $ws.Cells.Item(25, 3).Value = "浦银安盛精致生活混合"
    $ws.Cells.Item(25, 4).Value = "'2.09"
    $ws.Cells.Item(25, 5).Value = "'90.20"
    $ws.Cells.Item(25, 6).Value = "'2.84"
    $ws.Cells.Item(25, 7).Value = "'0.0594"
    $ws.Cells.Item(25, 8).Value = 6
    # row 26
    $ws.Cells.Item(26, 1).Value = 24
    $ws.Cells.Item(26, 2).Value = "'007307"
    $ws.Cells.Item(26, 3).Value = "华泰柏瑞基本面智选混合C"
    $ws.Cells.Item(26, 4).Value = "'1.03"
    $ws.Cells.Item(26, 5).Value = "'93.61"
    $ws.Cells.Item(26, 6).Value = "'4.16"
    $ws.Cells.Item(26, 7).Value = "'0.0428"
    $ws.Cells.Item(26, 8).Value = 10
    # row 27
    $ws.Cells.Item(27, 1).Value = 25
    $ws.Cells.Item(27, 2).Value = "'010346"
    $ws.Cells.Item(27, 3).Value = "华泰柏瑞成长智选混合C"
    $ws.Cells.Item(27, 4).Value = "'0.82"
    $ws.Cells.Item(27, 5).Value = "'93.47"
    $ws.Cells.Item(27, 6).Value = "'4.71"
    $ws.Cells.Item(27, 7).Value = "'0.0386"
    $ws.Cells.Item(27, 8).Value = 8
    # row 28
    $ws.Cells.Item(28, 1).Value = 26
    $ws.Cells.Item(28, 2).Value = "'004677"
    $ws.Cells.Item(28, 3).Value = "博时战略新兴产业混合"
    $ws.Cells.Item(28, 4).Value = "'0.41"
    $ws.Cells.Item(28, 5).Value = "'89.27"
    $ws.Cells.Item(28, 6).Value = "'4.84"
    $ws.Cells.Item(28, 7).Value = "'0.0198"
    $ws.Cells.Item(28, 8).Value = 7
    # row 29
    $ws.Cells.Item(29, 1).Value = 27
    $ws.Cells.Item(29, 2).Value = "'004223"
    $ws.Cells.Item(29, 3).Value = "金信多策略精选灵活配置混合"
    $ws.Cells.Item(29, 4).Value = "'0.36"
    $ws.Cells.Item(29, 5).Value = "'93.14"
    $ws.Cells.Item(29, 6).Value = "'5.27"
    $ws.Cells.Item(29, 7).Value = "'0.0190"
    $ws.Cells.Item(29, 8).Value = 5
    # row 30
    $ws.Cells.Item(30, 1).Value = 28
    $ws.Cells.Item(30, 2).Value = "'011462"
    $ws.Cells.Item(30, 3).Value = "东吴兴享成长混合C"
    $ws.Cells.Item(30, 4).Value = "'0.33"
    $ws.Cells.Item(30, 5).Value = "'80.15"
    $ws.Cells.Item(30, 6).Value = "'5.73"
    $ws.Cells.Item(30, 7).Value = "'0.0189"
    $ws.Cells.Item(30, 8).Value = 2
    # row 31
    $ws.Cells.Item(31, 1).Value = 29
    $ws.Cells.Item(31, 2).Value = "'004233"
    $ws.Cells.Item(31, 3).Value = "中欧盛世成长混合 (LOF) -C"
    $ws.Cells.Item(31, 4).Value = "'0.44"
    $ws.Cells.Item(31, 5).Value = "'85.98"
    $ws.Cells.Item(31, 6).Value = "'3.50"
    $ws.Cells.Item(31, 7).Value = "'0.0154"
    $ws.Cells.Item(31, 8).Value = 4
    # row 32
    $ws.Cells.Item(32, 1).Value = 30
    $ws.Cells.Item(32, 2).Value = "'001307"
    $ws.Cells.Item(32, 3).Value = "中欧永裕混合C"
    $ws.Cells.Item(32, 4).Value = "'0.35"
    $ws.Cells.Item(32, 5).Value = "'86.33"
    $ws.Cells.Item(32, 6).Value = "'3.52"
    $ws.Cells.Item(32, 7).Value = "'0.0123"
    $ws.Cells.Item(32, 8).Value = 4
    # row 33
    $ws.Cells.Item(33, 1).Value = 31
    $ws.Cells.Item(33, 2).Value = "'001731"
    $ws.Cells.Item(33, 3).Value = "广发百发大数据策略价值灵活配置混合A"
    $ws.Cells.Item(33, 4).Value = "'0.24"
    $ws.Cells.Item(33, 5).Value = "'88.87"
    $ws.Cells.Item(33, 6).Value = "'2.81"
    $ws.Cells.Item(33, 7).Value = "'0.0067"
    $ws.Cells.Item(33, 8).Value = 8
    # row 34
    $ws.Cells.Item(34, 1).Value = 32
    $ws.Cells.Item(34, 2).Value = "'001732"
    $ws.Cells.Item(34, 3).Value = "广发百发大数据策略价值灵活配置混合E"
    $ws.Cells.Item(34, 4).Value = "'0.24"
    $ws.Cells.Item(34, 5).Value = "'88.87"
    $ws.Cells.Item(34, 6).Value = "'2.81"
    $ws.Cells.Item(34, 7).Value = "'0.0067"
    $ws.Cells.Item(34, 8).Value = 8
    # row 35
    $ws.Cells.Item(35, 1).Value = 33
    $ws.Cells.Item(35, 2).Value = "'004917"
    $ws.Cells.Item(35, 3).Value = "中银证券祥瑞混合A"
    $ws.Cells.Item(35, 4).Value = "'0.10"
    $ws.Cells.Item(35, 5).Value = "'79.01"
    $ws.Cells.Item(35, 6).Value = "'2.06"
    $ws.Cells.Item(35, 7).Value = "'0.0021"
    $ws.Cells.Item(35, 8).Value = 9
    # row 36
    $ws.Cells.Item(36, 1).Value = 34
    $ws.Cells.Item(36, 2).Value = "'004918"
    $ws.Cells.Item(36, 3).Value = "中银证券祥瑞混合C"
    $ws.Cells.Item(36, 4).Value = "'0.07"
    $ws.Cells.Item(36, 5).Value = "'79.01"
    $ws.Cells.Item(36, 6).Value = "'2.06"
    $ws.Cells.Item(36, 7).Value = "'0.0014"
    $ws.Cells.Item(36, 8).Value = 9
    # row 37
    $ws.Cells.Item(37, 1).Value = 35
    $ws.Cells.Item(37, 2).Value = "'960031"
    $ws.Cells.Item(37, 3).Value = "浦银安盛价值成长混合H"
    $ws.Cells.Item(37, 4).Value = "'0.00"
    $ws.Cells.Item(37, 5).Value = "'88.87"
    $ws.Cells.Item(37, 6).Value = "'3.01"
    $ws.Cells.Item(37, 7).Value = 0
    $ws.Cells.Item(37, 8).Value = 3
    # row 38
    $ws.Cells.Item(38, 1).Value = 36
    $ws.Cells.Item(38, 2).Value = "'014011"
    $ws.Cells.Item(38, 3).Value = "浦银安盛价值成长混合C"
    $ws.Cells.Item(38, 4).Value = "'0.00"
    $ws.Cells.Item(38, 5).Value = "'88.87"
    $ws.Cells.Item(38, 6).Value = "'3.01"
    $ws.Cells.Item(38, 7).Value = 0
    $ws.Cells.Item(38, 8).Value = 3

# Re-use the existing header/index-column styling (bold, centered, bordered
# "s=2" style) from the "2021-Q4" sheet, which already has the identical
# layout, instead of constructing a brand-new style.
$styleSrc = $wb.Worksheets.Item("2021-Q4")
$styleSrc.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$ws.Range("A2:A38").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Prepend the 2022-Q1 row to the "总计" summary sheet
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows.Item(2).Insert()
$totalWs.Range("A2:D2").ClearFormats()

$totalWs.Cells.Item(2, 1).Value = 0
$totalWs.Cells.Item(2, 2).Value = "2022-Q1"
$totalWs.Cells.Item(2, 3).Value = 37
$totalWs.Cells.Item(2, 4).Value = 15.56

# Restore the index-column ("s=2") styling on the newly inserted row,
# copied from a sibling index cell.
$totalWs.Range("A3").Copy()
$totalWs.Range("A2").PasteSpecial(-4122)

# Renumber the index column (A) of the rows that shifted down by one.
for ($r = 3; $r -le 7; $r++) {
    $totalWs.Cells.Item($r, 1).Value = $r - 2
}
